$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 5939.875
$ws.Range("I19").Value = 673.3333
$ws.Range("K19").Value = 673.3333
$ws.Range("M19").Value = -498.3333

$ws.Range("H40").Value = 2649.3333
$ws.Range("I40").Value = 2098.75
$ws.Range("K40").Value = 2098.75
$ws.Range("M40").Value = -1923.75

$ws.Range("H113").Value = 102928
$ws.Range("I113").Value = 3314.8
$ws.Range("K113").Value = 3314.8
$ws.Range("M113").Value = -60.80000000000018

$ws.Range("H137").Value = 5823.8
$ws.Range("I137").Value = 2943.4375
$ws.Range("J137").Value = 10944.444
$ws.Range("K137").Value = 8830.3125
$ws.Range("L137").Value = 32833.33199999999
$ws.Range("M137").Value = -6280.3125
$ws.Range("N137").Value = -37933.33199999999

$ws.Range("H138").Value = 2813.725
$ws.Range("I138").Value = 1964.6316
$ws.Range("K138").Value = 5893.8948
$ws.Range("M138").Value = -753.8948

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 12826033
$ws.Range("I61").Value = 18521494
$ws.Range("K61").Value = 18521494
$ws.Range("M61").Value = -18521282

$ws.Range("H122").Value = 2001
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2001
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").Value = 6003
$ws.Range("N122").Value = -10903

$ws.Range("H136").Value = 12826033
$ws.Range("I136").Value = 18521494
$ws.Range("K136").Value = 55564482
$ws.Range("M136").Value = -55561932

$ws.Range("H141").Value = 96697.336
$ws.Range("J141").Value = 96696.5
$ws.Range("L141").Value = 96696.5
$ws.Range("N141").Value = -107056.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2646.1667
$ws.Range("I22").Value = 1694.8889
$ws.Range("J22").Value = 5500
$ws.Range("K22").Value = 1694.8889
$ws.Range("L22").Value = 5500
$ws.Range("M22").Value = -1521.8889
$ws.Range("N22").Value = -5846

$ws.Range("H82").Value = 36507.08
$ws.Range("I82").Value = 22791.7
$ws.Range("J82").Value = 82225
$ws.Range("K82").Value = 22791.7
$ws.Range("L82").Value = 82225
$ws.Range("M82").Value = -22408.7
$ws.Range("N82").Value = -82991

$ws.Range("H85").Value = 36507.08
$ws.Range("I85").Value = 22791.7
$ws.Range("J85").Value = 82225
$ws.Range("K85").Value = 22791.7
$ws.Range("L85").Value = 82225
$ws.Range("M85").Value = -21465.7
$ws.Range("N85").Value = -84877

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 43483380
$ws.Range("I31").Value = 111114070
$ws.Range("J31").Value = 6504
$ws.Range("K31").Value = 111114070
$ws.Range("L31").Value = 6504
$ws.Range("M31").Value = -111113775
$ws.Range("N31").Value = -7094

$ws.Range("H34").Value = 43483380
$ws.Range("I34").Value = 111114070
$ws.Range("J34").Value = 6504
$ws.Range("K34").Value = 111114070
$ws.Range("L34").Value = 6504
$ws.Range("M34").Value = -111113868
$ws.Range("N34").Value = -6908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2735
$ws.Range("I64").Value = 2060
$ws.Range("J64").Value = 2927.8572
$ws.Range("K64").Value = 6180
$ws.Range("L64").Value = 8783.571599999999
$ws.Range("M64").Value = -5910
$ws.Range("N64").Value = -9323.571599999999

$ws.Range("H67").Value = 2735
$ws.Range("I67").Value = 2060
$ws.Range("J67").Value = 2927.8572
$ws.Range("K67").Value = 6180
$ws.Range("L67").Value = 8783.571599999999
$ws.Range("M67").Value = -5244
$ws.Range("N67").Value = -10655.5716

$ws.Range("H86").Value = 330.7143
$ws.Range("J86").Value = 240
$ws.Range("L86").Value = 720
$ws.Range("N86").Value = -3092

$ws.Range("H89").Value = 330.7143
$ws.Range("J89").Value = 240
$ws.Range("L89").Value = 2160
$ws.Range("N89").Value = -14016

$ws.Range("H106").Value = 7506.5
$ws.Range("J106").Value = 9000
$ws.Range("L106").Value = 27000
$ws.Range("N106").Value = -28892

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 50000
$ws.Range("J69").Value = 50000
$ws.Range("L69").Value = 50000
$ws.Range("N69").Value = -51498

$ws.Range("H70").Value = 4666.5
$ws.Range("I70").Value = 4999.75
$ws.Range("K70").Value = 4999.75
$ws.Range("M70").Value = -4729.75

$ws.Range("H72").Value = 50000
$ws.Range("J72").Value = 50000
$ws.Range("L72").Value = 150000
$ws.Range("N72").Value = -157488

$ws.Range("H73").Value = 4666.5
$ws.Range("I73").Value = 4999.75
$ws.Range("K73").Value = 4999.75
$ws.Range("M73").Value = -4063.75

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").ClearContents()
$ws.Range("N87").Value = 0

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").ClearContents()
$ws.Range("N90").Value = 0

$ws.Range("H94").Value = 27499.5
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 27499.5
$ws.Range("K94").Value = 0
$ws.Range("L94").ClearContents()
$ws.Range("M94").Value = 27499.5
$ws.Range("N94").Value = -28851.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 963.3333
$ws.Range("I93").Value = 700
$ws.Range("K93").Value = 700
$ws.Range("M93").Value = 548

$ws.Range("H100").Value = 13161108
$ws.Range("I100").Value = 41668948
$ws.Range("K100").Value = 41668948
$ws.Range("M100").Value = -41668407

$ws.Range("H122").Value = 2783.7144
$ws.Range("I122").Value = 2783.7144
$ws.Range("K122").Value = 8351.143199999999
$ws.Range("M122").Value = -5901.143199999999

$ws.Range("H132").Value = 3034.2727
$ws.Range("I132").Value = 2159.6316
$ws.Range("J132").Value = 4989.353
$ws.Range("K132").Value = 6478.8948
$ws.Range("L132").Value = 14968.059
$ws.Range("M132").Value = -3948.8948
$ws.Range("N132").Value = -20028.059

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1101.5385
$ws.Range("I107").Value = 757.0968
$ws.Range("K107").Value = 2271.2904
$ws.Range("M107").Value = -351.2903999999999

$ws.Range("H122").Value = 3860.85
$ws.Range("I122").Value = 3961.2942
$ws.Range("J122").Value = 3291.6667
$ws.Range("K122").Value = 11883.8826
$ws.Range("L122").Value = 9875.000100000001
$ws.Range("M122").Value = -9433.882599999999
$ws.Range("N122").Value = -14775.0001
